# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.404.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.51%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.387.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.18%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.69%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.194"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.588"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.72"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.07%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.18%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "692.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.17%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.928.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.81%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.416.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.55%  "

# Row 17
$ws.Range("E17").Value = "  +1.79%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.363.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.73"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.38%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.904"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.75%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.65%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.99%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.62%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.60%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.70"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.70%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

# Row 31
$ws.Range("B31").Value = "dogwifhat"
$ws.Range("C31").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.75%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.39%  "

# Row 33
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "556.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.71%  "

# Row 34
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.02"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.36%  "

# Row 36
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.728.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.41%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.14%  "

# Row 38
$ws.Range("E38").Value = "  +8.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.78%  "

# Row 40
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0707"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.92%  "

# Row 41
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.40%  "

# Row 42
$ws.Range("B42").Value = "Fetch.AI"
$ws.Range("C42").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.67"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.57%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0419"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.37%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.130"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.50%  "

# Row 48
$ws.Range("E48").Value = "  -0.16%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.44"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.59%  "

# Row 51
$ws.Range("E51").Value = "  -2.98%  "
